$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates parsed directly from the target diff: (cell, newValue)
$updates = @(
    ,('D2', '36.602.06')
    ,('E2', '  +1.11%  ')
    ,('D3', '1.959.03')
    ,('E3', '  -0.04%  ')
    ,('E4', '  +0.05%  ')
    ,('D5', '244.53')
    ,('E5', '  +0.73%  ')
    ,('E6', '  +0.07%  ')
    ,('D7', '58.66')
    ,('E7', '  +1.10%  ')
    ,('E8', '  +0.02%  ')
    ,('E9', '  -0.95%  ')
    ,('D10', '56.36')
    ,('E10', '  -0.82%  ')
    ,('D11', '0.0864')
    ,('E11', '  +9.73%  ')
    ,('E12', '  +1.82%  ')
    ,('D13', '22.12')
    ,('E13', '  +2.31%  ')
    ,('E14', '  -1.55%  ')
    ,('D15', '2.245.18')
    ,('E15', '  -0.06%  ')
    ,('D16', '13.71')
    ,('E16', '  -1.59%  ')
    ,('D17', '5.25')
    ,('E17', '  -2.07%  ')
    ,('D18', '1.973.10')
    ,('E18', '  +0.45%  ')
    ,('D19', '36.514.24')
    ,('E19', '  +1.30%  ')
    ,('E20', '  +3.75%  ')
    ,('D21', '70.15')
    ,('E21', '  -1.18%  ')
    ,('D22', '230.19')
    ,('E22', '  -1.60%  ')
    ,('E23', '  -1.77%  ')
    ,('E24', '  -0.03%  ')
    ,('E25', '  -2.27%  ')
    ,('E26', '  +1.97%  ')
    ,('D27', '9.46')
    ,('E27', '  -1.40%  ')
    ,('D28', '162.37')
    ,('E28', '  +0.97%  ')
    ,('E29', '  +9.29%  ')
    ,('D30', '19.63')
    ,('E30', '  -0.70%  ')
    ,('D31', '0.119')
    ,('E31', '  -0.31%  ')
    ,('E32', '  +5.92%  ')
    ,('E33', '  -2.25%  ')
    ,('E34', '  +5.98%  ')
    ,('D35', '4.31')
    ,('E35', '  -1.75%  ')
    ,('D36', '6.42')
    ,('E36', '  +9.66%  ')
    ,('E37', '  +0.02%  ')
    ,('B38', 'WEMIXToken')
    ,('C38', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix')
    ,('D38', '1.78')
    ,('E38', '  -1.75%  ')
    ,('B39', 'LidoDAOToken')
    ,('C39', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo')
    ,('D39', '2.20')
    ,('E39', '  -2.91%  ')
    ,('D40', '3.05')
    ,('E40', '  +1.21%  ')
    ,('E41', '  +1.94%  ')
    ,('E42', '  +0.39%  ')
    ,('D43', '1.18')
    ,('E43', '  -2.35%  ')
    ,('E44', '  +0.32%  ')
    ,('D45', '16.22')
    ,('E45', '  +2.53%  ')
    ,('E46', '  -3.30%  ')
    ,('D47', '1.356.89')
    ,('E47', '  +1.46%  ')
    ,('D48', '88.71')
    ,('E48', '  -2.68%  ')
    ,('D49', '7.24')
    ,('E49', '  -3.47%  ')
    ,('D50', '2.82')
    ,('E50', '  +0.08%  ')
    ,('D51', '46.23')
    ,('E51', '  +5.21%  ')
)

foreach ($pair in $updates) {
    $cellRef = $pair[0]
    $newValue = $pair[1]
    $col = $cellRef -replace '[0-9]+$', ''
    $range = $ws.Range($cellRef)
    if ($col -eq 'D') {
        # Price column: values like "1.78" / "2.20" / "36.602.06" are strings in the
        # source data (multi-dot thousand separators etc.), not numbers. Force text
        # formatting before assignment so Excel does not coerce/round them, then drop
        # the temporary format so the cell keeps its original (unstyled) look.
        $range.NumberFormat = '@'
        $range.Value = $newValue
        $range.ClearFormats()
    } else {
        $range.Value = $newValue
    }
}
